$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 10:22:42"
$wsZh.Range("E3").Value = "2016-03-23 10:22:42"
$wsZh.Range("H2").Value = "2016-03-23 10:23:09"
$wsZh.Range("H3").Value = "2016-03-23 10:23:09"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 10:22:46"
$wsDe.Range("E3").Value = "2016-03-23 10:22:46"
$wsDe.Range("H2").Value = "2016-03-23 10:23:16"
$wsDe.Range("H3").Value = "2016-03-23 10:23:16"
